$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.864.69'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '1.562.53'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'205.94"
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  -1.96%  '
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("E10").Value = '  -1.37%  '
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("D12").Value = '1.785.06'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '1.572.75'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '26.877.37'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = "'61.22"
$ws.Range("E17").Value = '  -2.89%  '
$ws.Range("D18").Value = "'214.63"
$ws.Range("E18").Value = '  +1.30%  '
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").Value = "'9.16"
$ws.Range("E23").Value = '  -2.47%  '
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("D25").Value = "'153.84"
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("E26").Value = '  +2.08%  '
$ws.Range("D27").Value = "'14.91"
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  -0.87%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -3.10%  '
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("D33").Value = '1.401.77'
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("E34").Value = '  -0.73%  '
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("E36").Value = '  -1.16%  '
$ws.Range("D37").Value = "'0.922"
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").Value = "'0.525"
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").Value = "'0.996"
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = "'5.39"
$ws.Range("E43").Value = '  +3.61%  '
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("D45").Value = "'1.75"
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("D46").Value = "'63.12"
$ws.Range("E46").Value = '  -0.49%  '
$ws.Range("D47").Value = '1.697.97'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = "'86.33"
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("D49").Value = "'0.0506"
$ws.Range("E49").Value = '  +3.07%  '
$ws.Range("D50").Value = '0.0₇0980'
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("E51").Value = '  +0.55%  '
